$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Backlog")
$ws.Activate()

# Update Status column (L) for the sprint backlog rows
$ws.Range("L6").Value = "in progress"
$ws.Range("L7").Value = "in progress"
$ws.Range("L8").Value = "done"
$ws.Range("L9").Value = "in progress"
$ws.Range("L10").Value = "done"
$ws.Range("L11").Value = "done"
$ws.Range("L12").Value = "in progress"

# New actual-effort figures
$ws.Range("K8").Value = 0.33
$ws.Range("J9").Value = 8
$ws.Range("K10").Value = 0.33
$ws.Range("K11").Value = 0.33

# Move the active selection
$ws.Range("K6").Select() | Out-Null
